$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unidades")

$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 2

$ws.Range("F6").Select()
